# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Overview sheet: status text "Ready for handoff" -> "Handed back: in sync with en-US"
#  - zh-cn / de-de sheets: populate "Latest Target File" (hyperlinked) and
#    "Latest Handback File" / "Latest Handback DateTime" columns for both rows
#  - widen a few columns that now hold longer text

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# 1. Status text change (shared by Overview!E2:F3 and the Status column on
#    both language sheets) - update every occurrence so they stay in sync.
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus
$ws2.Range("C2").Value = $newStatus
$ws2.Range("C3").Value = $newStatus
$ws3.Range("C2").Value = $newStatus
$ws3.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn sheet (Item 2): fill in Latest Target File (I), Latest Handback
#    File (J) and Latest Handback DateTime (K) for both data rows.
# ---------------------------------------------------------------------------
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a671f92eca595d889a18c6177c2b97a8018a0610/e2e/a953fbbe-32ed-4438-bd54-f34e2800bf4d.md", [Type]::Missing, [Type]::Missing, "a953fbbe-32ed-4438-bd54-f34e2800bf4d.md") | Out-Null
$ws2.Range("J2").Value = "a953fbbe-32ed-4438-bd54-f34e2800bf4d.640354275b489af41758319e94de2ffc91410927.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-17 20:27:30"

$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a671f92eca595d889a18c6177c2b97a8018a0610/e2e/dfb9b19f-e2e6-471b-b44a-da0415d619bc.md", [Type]::Missing, [Type]::Missing, "dfb9b19f-e2e6-471b-b44a-da0415d619bc.md") | Out-Null
$ws2.Range("J3").Value = "dfb9b19f-e2e6-471b-b44a-da0415d619bc.73694072d94e2be2910f372da0b83d513673a408.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-08-17 20:27:30"

# ---------------------------------------------------------------------------
# 3. de-de sheet (Item 3): same three columns, different handback timestamp.
# ---------------------------------------------------------------------------
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a671f92eca595d889a18c6177c2b97a8018a0610/e2e/a953fbbe-32ed-4438-bd54-f34e2800bf4d.md", [Type]::Missing, [Type]::Missing, "a953fbbe-32ed-4438-bd54-f34e2800bf4d.md") | Out-Null
$ws3.Range("J2").Value = "a953fbbe-32ed-4438-bd54-f34e2800bf4d.640354275b489af41758319e94de2ffc91410927.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-17 20:27:37"

$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a671f92eca595d889a18c6177c2b97a8018a0610/e2e/dfb9b19f-e2e6-471b-b44a-da0415d619bc.md", [Type]::Missing, [Type]::Missing, "dfb9b19f-e2e6-471b-b44a-da0415d619bc.md") | Out-Null
$ws3.Range("J3").Value = "dfb9b19f-e2e6-471b-b44a-da0415d619bc.73694072d94e2be2910f372da0b83d513673a408.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-17 20:27:37"

# ---------------------------------------------------------------------------
# 4. Column widths: widen columns now holding the longer handback file names
#    / status text. ColumnWidth is in characters and snaps to a pixel grid,
#    so we pick the character width whose rounded pixel width lands on the
#    target column width (29.9777... -> 30, 18.65/21.71 -> 40).
# ---------------------------------------------------------------------------
$wide1 = 29.166666666666668   # -> stored width ~30 (was ~17.22)
$wide2 = 39.166666666666664   # -> stored width 40 (was ~18.65 / 21.71)

$ws1.Columns.Item(5).ColumnWidth = $wide1   # Overview!E
$ws1.Columns.Item(6).ColumnWidth = $wide1   # Overview!F

$ws2.Columns.Item(3).ColumnWidth = $wide1   # zh-cn!C (Status)
$ws2.Columns.Item(9).ColumnWidth = $wide2   # zh-cn!I (Latest Target File)
$ws2.Columns.Item(10).ColumnWidth = $wide2  # zh-cn!J (Latest Handback File)

$ws3.Columns.Item(3).ColumnWidth = $wide1   # de-de!C (Status)
$ws3.Columns.Item(9).ColumnWidth = $wide2   # de-de!I (Latest Target File)
$ws3.Columns.Item(10).ColumnWidth = $wide2  # de-de!J (Latest Handback File)
